# The post that used to live in row 759 ("「感情／感覚」") was removed from the
# sheet. Deleting the entire row shifts every subsequent row up by one
# (row 760 -> 759, row 761 -> 760, ... row 866 -> 865), which matches the
# reproduced diff, and Excel automatically shrinks the used range from
# A1:C866 to A1:C865.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(759).Delete()
